$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit reorders (permutes) the values of columns A,B,D,E,F,G,H,Q,R across
# rows 8-36, while leaving all other columns (C,I,J..P,S..AY) untouched.
# The target ("after") values for each touched cell are written directly below.

# Row 8
$ws.Cells.Item(8, 1).Value = 111756139
$ws.Cells.Item(8, 2).Value = 89405
$ws.Cells.Item(8, 4).Value = 'NT'
$ws.Cells.Item(8, 5).Value = 1202
$ws.Cells.Item(8, 6).Value = 'Ullticka'
$ws.Cells.Item(8, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(8, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(8, 17).Value = 453692.6056797595
$ws.Cells.Item(8, 18).Value = 7074032.491935454

# Row 10
$ws.Cells.Item(10, 1).Value = 111756140
$ws.Cells.Item(10, 2).Value = 89405
$ws.Cells.Item(10, 5).Value = 1202
$ws.Cells.Item(10, 6).Value = 'Ullticka'
$ws.Cells.Item(10, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(10, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(10, 17).Value = 453820.6239011836
$ws.Cells.Item(10, 18).Value = 7074037.242731699

# Row 12
$ws.Cells.Item(12, 1).Value = 111756166
$ws.Cells.Item(12, 2).Value = 77515
$ws.Cells.Item(12, 4).Value = 'NT'
$ws.Cells.Item(12, 5).Value = 6425
$ws.Cells.Item(12, 6).Value = 'Garnlav'
$ws.Cells.Item(12, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(12, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(12, 17).Value = 453981.6720900657
$ws.Cells.Item(12, 18).Value = 7073697.065866594

# Row 13
$ws.Cells.Item(13, 1).Value = 111756164
$ws.Cells.Item(13, 2).Value = 77515
$ws.Cells.Item(13, 5).Value = 6425
$ws.Cells.Item(13, 6).Value = 'Garnlav'
$ws.Cells.Item(13, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(13, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(13, 17).Value = 453971.0747186596
$ws.Cells.Item(13, 18).Value = 7073820.148138274

# Row 14
$ws.Cells.Item(14, 1).Value = 111756143
$ws.Cells.Item(14, 2).Value = 90087
$ws.Cells.Item(14, 4).Value = 'LC'
$ws.Cells.Item(14, 5).Value = 3298
$ws.Cells.Item(14, 6).Value = 'Trådticka'
$ws.Cells.Item(14, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(14, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(14, 17).Value = 453950.9091414157
$ws.Cells.Item(14, 18).Value = 7073591.829928016

# Row 15
$ws.Cells.Item(15, 1).Value = 111756168
$ws.Cells.Item(15, 2).Value = 77515
$ws.Cells.Item(15, 5).Value = 6425
$ws.Cells.Item(15, 6).Value = 'Garnlav'
$ws.Cells.Item(15, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(15, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(15, 17).Value = 453958.9423245317
$ws.Cells.Item(15, 18).Value = 7073596.134472342

# Row 16
$ws.Cells.Item(16, 1).Value = 111756171
$ws.Cells.Item(16, 2).Value = 88899
$ws.Cells.Item(16, 4).Value = 'NT'
$ws.Cells.Item(16, 5).Value = 3286
$ws.Cells.Item(16, 6).Value = 'Flattoppad klubbsvamp'
$ws.Cells.Item(16, 7).Value = 'Clavariadelphus truncatus'
$ws.Cells.Item(16, 8).Value = '(Quél.) Donk'
$ws.Cells.Item(16, 17).Value = 453750.6060291855
$ws.Cells.Item(16, 18).Value = 7073942.323881648

# Row 17
$ws.Cells.Item(17, 1).Value = 111756156
$ws.Cells.Item(17, 2).Value = 89423
$ws.Cells.Item(17, 5).Value = 5432
$ws.Cells.Item(17, 6).Value = 'Granticka'
$ws.Cells.Item(17, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(17, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(17, 17).Value = 453978.4965374623
$ws.Cells.Item(17, 18).Value = 7073812.964766338

# Row 18
$ws.Cells.Item(18, 1).Value = 111756170
$ws.Cells.Item(18, 2).Value = 96265
$ws.Cells.Item(18, 4).Value = 'LC'
$ws.Cells.Item(18, 5).Value = 219790
$ws.Cells.Item(18, 6).Value = 'Fläcknycklar'
$ws.Cells.Item(18, 7).Value = 'Dactylorhiza maculata'
$ws.Cells.Item(18, 8).Value = '(L.) Soó'
$ws.Cells.Item(18, 17).Value = 453738.5427278728
$ws.Cells.Item(18, 18).Value = 7073724.066700204

# Row 19
$ws.Cells.Item(19, 1).Value = 111756155
$ws.Cells.Item(19, 2).Value = 89423
$ws.Cells.Item(19, 5).Value = 5432
$ws.Cells.Item(19, 6).Value = 'Granticka'
$ws.Cells.Item(19, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(19, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(19, 17).Value = 453863.4009631127
$ws.Cells.Item(19, 18).Value = 7073965.428905412

# Row 20
$ws.Cells.Item(20, 1).Value = 111756162
$ws.Cells.Item(20, 2).Value = 77515
$ws.Cells.Item(20, 5).Value = 6425
$ws.Cells.Item(20, 6).Value = 'Garnlav'
$ws.Cells.Item(20, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(20, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(20, 17).Value = 453922.6243923472
$ws.Cells.Item(20, 18).Value = 7073958.370937477

# Row 21
$ws.Cells.Item(21, 1).Value = 111756148
$ws.Cells.Item(21, 2).Value = 96266
$ws.Cells.Item(21, 4).Value = 'LC'
$ws.Cells.Item(21, 5).Value = 223591
$ws.Cells.Item(21, 6).Value = 'Skogsnycklar'
$ws.Cells.Item(21, 7).Value = 'Dactylorhiza maculata subsp. fuchsii'
$ws.Cells.Item(21, 8).Value = '(Druce) Hyl.'
$ws.Cells.Item(21, 17).Value = 453747.0542679164
$ws.Cells.Item(21, 18).Value = 7073851.289854143

# Row 22
$ws.Cells.Item(22, 1).Value = 111756147
$ws.Cells.Item(22, 2).Value = 89425
$ws.Cells.Item(22, 4).Value = 'NT'
$ws.Cells.Item(22, 5).Value = 5442
$ws.Cells.Item(22, 6).Value = 'Tallticka'
$ws.Cells.Item(22, 7).Value = 'Porodaedalea pini'
$ws.Cells.Item(22, 8).Value = '(Brot.) Murrill'
$ws.Cells.Item(22, 17).Value = 453989.3915585176
$ws.Cells.Item(22, 18).Value = 7073710.21875874

# Row 23
$ws.Cells.Item(23, 1).Value = 111756153
$ws.Cells.Item(23, 2).Value = 96674
$ws.Cells.Item(23, 4).Value = 'LC'
$ws.Cells.Item(23, 5).Value = 219880
$ws.Cells.Item(23, 6).Value = 'Kransrams'
$ws.Cells.Item(23, 7).Value = 'Polygonatum verticillatum'
$ws.Cells.Item(23, 8).Value = '(L.) All.'
$ws.Cells.Item(23, 17).Value = 453707.5163784204
$ws.Cells.Item(23, 18).Value = 7073721.869806641

# Row 25
$ws.Cells.Item(25, 1).Value = 111756151
$ws.Cells.Item(25, 2).Value = 95532
$ws.Cells.Item(25, 4).Value = 'LC'
$ws.Cells.Item(25, 5).Value = 221945
$ws.Cells.Item(25, 6).Value = 'Revlummer'
$ws.Cells.Item(25, 7).Value = 'Lycopodium annotinum'
$ws.Cells.Item(25, 8).Value = 'L.'
$ws.Cells.Item(25, 17).Value = 453609.4901279925
$ws.Cells.Item(25, 18).Value = 7074130.545069677

# Row 26
$ws.Cells.Item(26, 1).Value = 111756157
$ws.Cells.Item(26, 2).Value = 89423
$ws.Cells.Item(26, 4).Value = 'NT'
$ws.Cells.Item(26, 5).Value = 5432
$ws.Cells.Item(26, 6).Value = 'Granticka'
$ws.Cells.Item(26, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(26, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(26, 17).Value = 453981.5111392652
$ws.Cells.Item(26, 18).Value = 7073807.172376178

# Row 27
$ws.Cells.Item(27, 1).Value = 111756165
$ws.Cells.Item(27, 17).Value = 453984.2379404157
$ws.Cells.Item(27, 18).Value = 7073751.417626478

# Row 28
$ws.Cells.Item(28, 1).Value = 111756150
$ws.Cells.Item(28, 2).Value = 95532
$ws.Cells.Item(28, 5).Value = 221945
$ws.Cells.Item(28, 6).Value = 'Revlummer'
$ws.Cells.Item(28, 7).Value = 'Lycopodium annotinum'
$ws.Cells.Item(28, 8).Value = 'L.'
$ws.Cells.Item(28, 17).Value = 453976.2702886119
$ws.Cells.Item(28, 18).Value = 7073812.112971266

# Row 29
$ws.Cells.Item(29, 1).Value = 111756154
$ws.Cells.Item(29, 2).Value = 96674
$ws.Cells.Item(29, 4).Value = 'LC'
$ws.Cells.Item(29, 5).Value = 219880
$ws.Cells.Item(29, 6).Value = 'Kransrams'
$ws.Cells.Item(29, 7).Value = 'Polygonatum verticillatum'
$ws.Cells.Item(29, 8).Value = '(L.) All.'
$ws.Cells.Item(29, 17).Value = 453614.9183513908
$ws.Cells.Item(29, 18).Value = 7074108.35826167

# Row 30
$ws.Cells.Item(30, 1).Value = 111756142
$ws.Cells.Item(30, 2).Value = 90087
$ws.Cells.Item(30, 4).Value = 'LC'
$ws.Cells.Item(30, 5).Value = 3298
$ws.Cells.Item(30, 6).Value = 'Trådticka'
$ws.Cells.Item(30, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(30, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(30, 17).Value = 454002.5104495964
$ws.Cells.Item(30, 18).Value = 7073638.391199326

# Row 31
$ws.Cells.Item(31, 1).Value = 111756160
$ws.Cells.Item(31, 2).Value = 77515
$ws.Cells.Item(31, 5).Value = 6425
$ws.Cells.Item(31, 6).Value = 'Garnlav'
$ws.Cells.Item(31, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(31, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(31, 17).Value = 453815.5156181521
$ws.Cells.Item(31, 18).Value = 7073870.182023689

# Row 32
$ws.Cells.Item(32, 1).Value = 111756163
$ws.Cells.Item(32, 17).Value = 453955.6479769219
$ws.Cells.Item(32, 18).Value = 7073945.9492877

# Row 33
$ws.Cells.Item(33, 1).Value = 111756141
$ws.Cells.Item(33, 2).Value = 89405
$ws.Cells.Item(33, 4).Value = 'NT'
$ws.Cells.Item(33, 5).Value = 1202
$ws.Cells.Item(33, 6).Value = 'Ullticka'
$ws.Cells.Item(33, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(33, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(33, 17).Value = 453610.1793069927
$ws.Cells.Item(33, 18).Value = 7074087.205471905

# Row 34
$ws.Cells.Item(34, 1).Value = 111756169
$ws.Cells.Item(34, 2).Value = 77515
$ws.Cells.Item(34, 4).Value = 'NT'
$ws.Cells.Item(34, 5).Value = 6425
$ws.Cells.Item(34, 6).Value = 'Garnlav'
$ws.Cells.Item(34, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(34, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(34, 17).Value = 453910.2023238647
$ws.Cells.Item(34, 18).Value = 7073654.334338664

# Row 35
$ws.Cells.Item(35, 1).Value = 111756161
$ws.Cells.Item(35, 2).Value = 77515
$ws.Cells.Item(35, 5).Value = 6425
$ws.Cells.Item(35, 6).Value = 'Garnlav'
$ws.Cells.Item(35, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(35, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(35, 17).Value = 453723.2573215028
$ws.Cells.Item(35, 18).Value = 7074069.623294062

# Row 36
$ws.Cells.Item(36, 1).Value = 111756158
$ws.Cells.Item(36, 2).Value = 89423
$ws.Cells.Item(36, 5).Value = 5432
$ws.Cells.Item(36, 6).Value = 'Granticka'
$ws.Cells.Item(36, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(36, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(36, 17).Value = 454002.8592168373
$ws.Cells.Item(36, 18).Value = 7073783.424762985

